$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = "oi_short"
$ws.Range("B4").Value  = "pool_balance_btc"
$ws.Range("B5").Value  = "cum_apy_providers"
$ws.Range("B6").Value  = "volume_eth"
$ws.Range("B7").Value  = "number_of_liquidations"
$ws.Range("B8").Value  = "volume_sol"
$ws.Range("B9").Value  = "treasury_balance"
$ws.Range("B10").Value = "pool_balance_usdc"
$ws.Range("B11").Value = "min_pnl_traders"
$ws.Range("B12").Value = "number_of_traders"
$ws.Range("B13").Value = "pool_balance_usdT"
$ws.Range("B14").Value = "fees_collected"
$ws.Range("B15").Value = "cum_pnl_traders"
$ws.Range("B16").Value = "volume_btc"
$ws.Range("B17").Value = "pool_balance_eth"
$ws.Range("B18").Value = "oi_long"
$ws.Range("B19").Value = "pool_balance_sol"
$ws.Range("B20").Value = "max_pnl_traders"
$ws.Range("B21").Value = "number_of_liquidity_providers"
